$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value2 = 0.296
$ws.Range("C4").Value2 = 0.051
$ws.Range("E4").Value2 = 0.151
$ws.Range("H4").Value2 = 0.189
$ws.Range("J4").Value2 = 0.108
$ws.Range("K4").Value2 = 0.342
$ws.Range("L4").Value2 = 0.102
$ws.Range("M4").Value2 = 0.32
$ws.Range("N4").Value2 = 0.271
$ws.Range("O4").Value2 = 0.02
$ws.Range("P4").Value2 = 0.143
$ws.Range("Q4").Value2 = 0.513
$ws.Range("R4").Value2 = 0.217
$ws.Range("S4").Value2 = 0.466
$ws.Range("T4").Value2 = 0.283
$ws.Range("W4").Value2 = 0.244
$ws.Range("Y4").Value2 = 0.208
$ws.Range("Z4").Value2 = 0.451
$ws.Range("AA4").Value2 = 0.133
$ws.Range("AB4").Value2 = 0.364
$ws.Range("AC4").Value2 = 0.126
$ws.Range("AE4").Value2 = 0.078
$ws.Range("AF4").Value2 = 0.737
$ws.Range("AG4").Value2 = 0.094
$ws.Range("AH4").Value2 = 0.307
$ws.Range("AI4").Value2 = 0.658
$ws.Range("AJ4").Value2 = 0.172
$ws.Range("AK4").Value2 = 0.415
$ws.Range("AL4").Value2 = 0.703
$ws.Range("AN4").Value2 = 0.339
$ws.Range("AO4").Value2 = 0.699

# Row 5
$ws.Range("B5").Value2 = 0.8159999999999999
$ws.Range("C5").Value2 = 0.15
$ws.Range("D5").Value2 = 0.388
$ws.Range("E5").Value2 = 0.658
$ws.Range("F5").Value2 = 0.225
$ws.Range("G5").Value2 = 0.474
$ws.Range("H5").Value2 = 0.8159999999999999
$ws.Range("I5").Value2 = 0.15
$ws.Range("J5").Value2 = 0.388
$ws.Range("K5").Value2 = 0.658
$ws.Range("L5").Value2 = 0.225
$ws.Range("M5").Value2 = 0.474
$ws.Range("N5").Value2 = 0.842
$ws.Range("O5").Value2 = 0.133
$ws.Range("P5").Value2 = 0.365
$ws.Range("Q5").Value2 = 0.579
$ws.Range("R5").Value2 = 0.244
$ws.Range("S5").Value2 = 0.494
$ws.Range("T5").Value2 = 0.579
$ws.Range("U5").Value2 = 0.244
$ws.Range("V5").Value2 = 0.494
$ws.Range("W5").Value2 = 0.737
$ws.Range("X5").Value2 = 0.194
$ws.Range("Y5").Value2 = 0.44
$ws.Range("Z5").Value2 = 0.8159999999999999
$ws.Range("AA5").Value2 = 0.15
$ws.Range("AB5").Value2 = 0.388
$ws.Range("AC5").Value2 = 0.763
$ws.Range("AD5").Value2 = 0.181
$ws.Range("AE5").Value2 = 0.425
$ws.Range("AF5").Value2 = 0.974
$ws.Range("AH5").Value2 = 0.16
$ws.Range("AI5").Value2 = 0.763
$ws.Range("AJ5").Value2 = 0.181
$ws.Range("AK5").Value2 = 0.425
$ws.Range("AL5").Value2 = 0.921
$ws.Range("AM5").Value2 = 0.073
$ws.Range("AN5").Value2 = 0.27
$ws.Range("AO5").Value2 = 0.886

# Row 6
$ws.Range("B6").Value2 = 0.434
$ws.Range("E6").Value2 = 0.246
$ws.Range("H6").Value2 = 0.307
$ws.Range("K6").Value2 = 0.45
$ws.Range("N6").Value2 = 0.41
$ws.Range("Q6").Value2 = 0.544
$ws.Range("T6").Value2 = 0.38
$ws.Range("W6").Value2 = 0.367
$ws.Range("Z6").Value2 = 0.581
$ws.Range("AC6").Value2 = 0.216
$ws.Range("AF6").Value2 = 0.839
$ws.Range("AI6").Value2 = 0.707
$ws.Range("AL6").Value2 = 0.797
$ws.Range("AO6").Value2 = 0.781

# Row 7
$ws.Range("B7").Value2 = 0.604
$ws.Range("E7").Value2 = 0.394
$ws.Range("H7").Value2 = 0.491
$ws.Range("K7").Value2 = 0.555
$ws.Range("N7").Value2 = 0.592
$ws.Range("Q7").Value2 = 0.5639999999999999
$ws.Range("T7").Value2 = 0.479
$ws.Range("W7").Value2 = 0.525
$ws.Range("Z7").Value2 = 0.702
$ws.Range("AC7").Value2 = 0.379
$ws.Range("AF7").Value2 = 0.915
$ws.Range("AI7").Value2 = 0.739
$ws.Range("AL7").Value2 = 0.867
$ws.Range("AO7").Value2 = 0.84

# Row 8
$ws.Range("B8").Value2 = 0.751
$ws.Range("C8").Value2 = 0.15
$ws.Range("D8").Value2 = 0.387
$ws.Range("E8").Value2 = 0.548
$ws.Range("F8").Value2 = 0.191
$ws.Range("G8").Value2 = 0.437
$ws.Range("H8").Value2 = 0.703
$ws.Range("I8").Value2 = 0.152
$ws.Range("J8").Value2 = 0.39
$ws.Range("K8").Value2 = 0.586
$ws.Range("L8").Value2 = 0.204
$ws.Range("M8").Value2 = 0.452
$ws.Range("N8").Value2 = 0.749
$ws.Range("O8").Value2 = 0.137
$ws.Range("P8").Value2 = 0.371
$ws.Range("Q8").Value2 = 0.55
$ws.Range("R8").Value2 = 0.229
$ws.Range("S8").Value2 = 0.479
$ws.Range("T8").Value2 = 0.501
$ws.Range("U8").Value2 = 0.206
$ws.Range("V8").Value2 = 0.454
$ws.Range("W8").Value2 = 0.662
$ws.Range("X8").Value2 = 0.18
$ws.Range("Y8").Value2 = 0.424
$ws.Range("Z8").Value2 = 0.751
$ws.Range("AA8").Value2 = 0.15
$ws.Range("AB8").Value2 = 0.387
$ws.Range("AC8").Value2 = 0.655
$ws.Range("AD8").Value2 = 0.177
$ws.Range("AE8").Value2 = 0.42
$ws.Range("AF8").Value2 = 0.893
$ws.Range("AG8").Value2 = 0.046
$ws.Range("AH8").Value2 = 0.215
$ws.Range("AI8").Value2 = 0.753
$ws.Range("AJ8").Value2 = 0.18
$ws.Range("AK8").Value2 = 0.424
$ws.Range("AL8").Value2 = 0.892
$ws.Range("AM8").Value2 = 0.078
$ws.Range("AN8").Value2 = 0.279
$ws.Range("AO8").Value2 = 0.846

# Row 9
$ws.Range("B9").Value2 = 0.658
$ws.Range("C9").Value2 = 0.225
$ws.Range("D9").Value2 = 0.474
$ws.Range("E9").Value2 = 0.421
$ws.Range("F9").Value2 = 0.244
$ws.Range("G9").Value2 = 0.494
$ws.Range("H9").Value2 = 0.579
$ws.Range("I9").Value2 = 0.244
$ws.Range("J9").Value2 = 0.494
$ws.Range("K9").Value2 = 0.5
$ws.Range("N9").Value2 = 0.632
$ws.Range("O9").Value2 = 0.233
$ws.Range("P9").Value2 = 0.482
$ws.Range("Q9").Value2 = 0.5
$ws.Range("T9").Value2 = 0.395
$ws.Range("U9").Value2 = 0.239
$ws.Range("V9").Value2 = 0.489
$ws.Range("W9").Value2 = 0.553
$ws.Range("X9").Value2 = 0.247
$ws.Range("Y9").Value2 = 0.497
$ws.Range("Z9").Value2 = 0.658
$ws.Range("AA9").Value2 = 0.225
$ws.Range("AB9").Value2 = 0.474
$ws.Range("AC9").Value2 = 0.553
$ws.Range("AD9").Value2 = 0.247
$ws.Range("AE9").Value2 = 0.497
$ws.Range("AF9").Value2 = 0.763
$ws.Range("AG9").Value2 = 0.181
$ws.Range("AH9").Value2 = 0.425
$ws.Range("AI9").Value2 = 0.737
$ws.Range("AJ9").Value2 = 0.194
$ws.Range("AK9").Value2 = 0.44
$ws.Range("AL9").Value2 = 0.842
$ws.Range("AM9").Value2 = 0.133
$ws.Range("AN9").Value2 = 0.365
$ws.Range("AO9").Value2 = 0.781

# Row 10
$ws.Range("B10").Value2 = 0.8159999999999999
$ws.Range("C10").Value2 = 0.15
$ws.Range("D10").Value2 = 0.388
$ws.Range("E10").Value2 = 0.579
$ws.Range("F10").Value2 = 0.244
$ws.Range("G10").Value2 = 0.494
$ws.Range("H10").Value2 = 0.737
$ws.Range("I10").Value2 = 0.194
$ws.Range("J10").Value2 = 0.44
$ws.Range("K10").Value2 = 0.658
$ws.Range("L10").Value2 = 0.225
$ws.Range("M10").Value2 = 0.474
$ws.Range("N10").Value2 = 0.8159999999999999
$ws.Range("O10").Value2 = 0.15
$ws.Range("P10").Value2 = 0.388
$ws.Range("Q10").Value2 = 0.579
$ws.Range("R10").Value2 = 0.244
$ws.Range("S10").Value2 = 0.494
$ws.Range("T10").Value2 = 0.579
$ws.Range("U10").Value2 = 0.244
$ws.Range("V10").Value2 = 0.494
$ws.Range("W10").Value2 = 0.737
$ws.Range("X10").Value2 = 0.194
$ws.Range("Y10").Value2 = 0.44
$ws.Range("Z10").Value2 = 0.8159999999999999
$ws.Range("AA10").Value2 = 0.15
$ws.Range("AB10").Value2 = 0.388
$ws.Range("AC10").Value2 = 0.658
$ws.Range("AD10").Value2 = 0.225
$ws.Range("AE10").Value2 = 0.474
$ws.Range("AF10").Value2 = 0.974
$ws.Range("AH10").Value2 = 0.16
$ws.Range("AI10").Value2 = 0.763
$ws.Range("AJ10").Value2 = 0.181
$ws.Range("AK10").Value2 = 0.425
$ws.Range("AL10").Value2 = 0.921
$ws.Range("AM10").Value2 = 0.073
$ws.Range("AN10").Value2 = 0.27
$ws.Range("AO10").Value2 = 0.886

# Row 11
$ws.Range("B11").Value2 = 0.8159999999999999
$ws.Range("C11").Value2 = 0.15
$ws.Range("D11").Value2 = 0.388
$ws.Range("E11").Value2 = 0.658
$ws.Range("F11").Value2 = 0.225
$ws.Range("G11").Value2 = 0.474
$ws.Range("H11").Value2 = 0.8159999999999999
$ws.Range("I11").Value2 = 0.15
$ws.Range("J11").Value2 = 0.388
$ws.Range("K11").Value2 = 0.658
$ws.Range("L11").Value2 = 0.225
$ws.Range("M11").Value2 = 0.474
$ws.Range("N11").Value2 = 0.842
$ws.Range("O11").Value2 = 0.133
$ws.Range("P11").Value2 = 0.365
$ws.Range("Q11").Value2 = 0.579
$ws.Range("R11").Value2 = 0.244
$ws.Range("S11").Value2 = 0.494
$ws.Range("T11").Value2 = 0.579
$ws.Range("U11").Value2 = 0.244
$ws.Range("V11").Value2 = 0.494
$ws.Range("W11").Value2 = 0.737
$ws.Range("X11").Value2 = 0.194
$ws.Range("Y11").Value2 = 0.44
$ws.Range("Z11").Value2 = 0.8159999999999999
$ws.Range("AA11").Value2 = 0.15
$ws.Range("AB11").Value2 = 0.388
$ws.Range("AC11").Value2 = 0.711
$ws.Range("AD11").Value2 = 0.206
$ws.Range("AE11").Value2 = 0.454
$ws.Range("AF11").Value2 = 0.974
$ws.Range("AH11").Value2 = 0.16
$ws.Range("AI11").Value2 = 0.763
$ws.Range("AJ11").Value2 = 0.181
$ws.Range("AK11").Value2 = 0.425
$ws.Range("AL11").Value2 = 0.921
$ws.Range("AM11").Value2 = 0.073
$ws.Range("AN11").Value2 = 0.27
$ws.Range("AO11").Value2 = 0.886

# Row 12
$ws.Range("B12").Value2 = 1.258
$ws.Range("C12").Value2 = 0.32
$ws.Range("D12").Value2 = 0.5659999999999999
$ws.Range("E12").Value2 = 1.68
$ws.Range("F12").Value2 = 1.098
$ws.Range("G12").Value2 = 1.048
$ws.Range("H12").Value2 = 1.613
$ws.Range("I12").Value2 = 1.334
$ws.Range("J12").Value2 = 1.155
$ws.Range("K12").Value2 = 1.4
$ws.Range("L12").Value2 = 0.5600000000000001
$ws.Range("M12").Value2 = 0.748
$ws.Range("N12").Value2 = 1.406
$ws.Range("O12").Value2 = 0.616
$ws.Range("P12").Value2 = 0.785
$ws.Range("Z12").Value2 = 1.258
$ws.Range("AA12").Value2 = 0.32
$ws.Range("AB12").Value2 = 0.5659999999999999
$ws.Range("AC12").Value2 = 1.793
$ws.Range("AD12").Value2 = 2.44
$ws.Range("AE12").Value2 = 1.562
$ws.Range("AF12").Value2 = 1.243
$ws.Range("AG12").Value2 = 0.238
$ws.Range("AH12").Value2 = 0.488
$ws.Range("AI12").Value2 = 1.034
$ws.Range("AJ12").Value2 = 0.033
$ws.Range("AK12").Value2 = 0.182
$ws.Range("AL12").Value2 = 1.086
$ws.Range("AM12").Value2 = 0.078
$ws.Range("AN12").Value2 = 0.28
$ws.Range("AO12").Value2 = 1.121

# Row 13
$ws.Range("B13").Value2 = 3.474
$ws.Range("C13").Value2 = 1.46
$ws.Range("D13").Value2 = 1.208
$ws.Range("E13").Value2 = 4.594
$ws.Range("F13").Value2 = 0.429
$ws.Range("G13").Value2 = 0.655
$ws.Range("H13").Value2 = 4.611
$ws.Range("I13").Value2 = 0.627
$ws.Range("J13").Value2 = 0.792
$ws.Range("K13").Value2 = 2.265
$ws.Range("L13").Value2 = 0.606
$ws.Range("M13").Value2 = 0.779
$ws.Range("N13").Value2 = 3.263
$ws.Range("O13").Value2 = 0.72
$ws.Range("P13").Value2 = 0.849
$ws.Range("Z13").Value2 = 2.514
$ws.Range("AA13").Value2 = 2.878
$ws.Range("AB13").Value2 = 1.697
$ws.Range("AC13").Value2 = 6.378
$ws.Range("AD13").Value2 = 2.181
$ws.Range("AE13").Value2 = 1.477
$ws.Range("AF13").Value2 = 1.605
$ws.Range("AG13").Value2 = 0.713
$ws.Range("AH13").Value2 = 0.844
$ws.Range("AI13").Value2 = 1.289
$ws.Range("AJ13").Value2 = 0.364
$ws.Range("AK13").Value2 = 0.603
$ws.Range("AL13").Value2 = 1.579
$ws.Range("AM13").Value2 = 0.717
$ws.Range("AN13").Value2 = 0.847
$ws.Range("AO13").Value2 = 1.491
